$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("KPIs")
$ws.Rows.Item(13).RowHeight = 28.35
